$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = "Иван Петров"
$ws.Range("C2").Value = "01.08.2020 0:00:00"
$ws.Range("E2").Value = "Выхлопная труба, Выхлопная труба, Выхлопная труба, Кузов, "

# --- Row 3 ---
$ws.Range("B3").Value = "Андрей Кругов"
$ws.Range("C3").Value = "01.08.2020 0:00:00"
$ws.Range("E3").Value = "Выхлопная труба, "

# D2 / D3 need the numeric-looking text "15" / "10" stored as TEXT (matching
# the source data), not auto-detected as numbers. Build each value through a
# helper formula (forces a text result) and paste-special just the value, so
# the target cell keeps its plain default style/format (no NumberFormat
# change lingers behind on the cell or in the workbook's style table).
$helper = $ws.Range("Z1")

$helper.Formula = '=TEXT(15,"0")'
$helper.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$helper.Formula = '=TEXT(10,"0")'
$helper.Copy()
$ws.Range("D3").PasteSpecial(-4163)

$helper.ClearContents()
